$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update loading percent results for the 380 kV case (res_line/loading_percent)

$ws.Range("B2").Value = 21.33417609464704
$ws.Range("C2").Value = 13.73195939775041
$ws.Range("D2").Value = 5.42443164944441
$ws.Range("E2").Value = 11.71371822519192
$ws.Range("F2").Value = 56.73657998384415
$ws.Range("J2").Value = 10.62657631756334

$ws.Range("B3").Value = 21.01014982325743
$ws.Range("C3").Value = 13.44516397771525
$ws.Range("D3").Value = 5.441226985526482
$ws.Range("E3").Value = 11.74800703038231
$ws.Range("F3").Value = 55.82059686701005
$ws.Range("J3").Value = 10.62669444130072

$ws.Range("B4").Value = 20.81808330194159
$ws.Range("C4").Value = 13.27317944290382
$ws.Range("D4").Value = 5.452785114294226
$ws.Range("E4").Value = 11.77214752248055
$ws.Range("F4").Value = 55.26172774199275
$ws.Range("J4").Value = 10.6287120180379

$ws.Range("B5").Value = 20.74166084315012
$ws.Range("C5").Value = 13.20425318476527
$ws.Range("D5").Value = 5.457805813967387
$ws.Range("E5").Value = 11.78275943283505
$ws.Range("F5").Value = 55.03509706198321
$ws.Range("J5").Value = 10.63002224102857

$ws.Range("B6").Value = 20.72908571583854
$ws.Range("C6").Value = 13.19288166837758
$ws.Range("D6").Value = 5.458658185592003
$ws.Range("E6").Value = 11.78456825354816
$ws.Range("F6").Value = 54.99753884496895
$ws.Range("J6").Value = 10.63026924481121

$ws.Range("B7").Value = 20.81704501867739
$ws.Range("C7").Value = 13.27224502208256
$ws.Range("D7").Value = 5.452851570715714
$ws.Range("E7").Value = 11.77228750506551
$ws.Range("F7").Value = 55.25866652454599
$ws.Range("J7").Value = 10.62872771360531

$ws.Range("B8").Value = 21.22109236939054
$ws.Range("C8").Value = 13.63229510556415
$ws.Range("D8").Value = 5.429962578166369
$ws.Range("E8").Value = 11.72489924473359
$ws.Range("F8").Value = 56.42014548838525
$ws.Range("J8").Value = 10.62621276812223

$ws.Range("B9").Value = 22.06286437968508
$ws.Range("C9").Value = 14.36542828665396
$ws.Range("D9").Value = 5.395077532701703
$ws.Range("E9").Value = 11.65654190538991
$ws.Range("F9").Value = 58.71623086377641
$ws.Range("J9").Value = 10.6367552520322

$ws.Range("B10").Value = 22.70433347311664
$ws.Range("C10").Value = 14.91328305908673
$ws.Range("D10").Value = 5.375703544588938
$ws.Range("E10").Value = 11.62139843535167
$ws.Range("F10").Value = 60.4014863376362
$ws.Range("J10").Value = 10.65398212094585

$ws.Range("B11").Value = 22.99963746556543
$ws.Range("C11").Value = 15.16304597305385
$ws.Range("D11").Value = 5.368282667360548
$ws.Range("E11").Value = 11.60870379205063
$ws.Range("F11").Value = 61.1651198114189
$ws.Range("J11").Value = 10.66388420803146

$ws.Range("B12").Value = 23.1118400462183
$ws.Range("C12").Value = 15.25758697405174
$ws.Range("D12").Value = 5.365675639067001
$ws.Range("E12").Value = 11.60437138719287
$ws.Range("F12").Value = 61.45364111745116
$ws.Range("J12").Value = 10.66793099940383

$ws.Range("B13").Value = 23.08766039212401
$ws.Range("C13").Value = 15.23722939185923
$ws.Range("D13").Value = 5.366228030820611
$ws.Range("E13").Value = 11.60528331210944
$ws.Range("F13").Value = 61.39153521181683
$ws.Range("J13").Value = 10.66704623857275

$ws.Range("B14").Value = 23.00886158348295
$ws.Range("C14").Value = 15.1708253408007
$ws.Range("D14").Value = 5.368064098793146
$ws.Range("E14").Value = 11.60833783995911
$ws.Range("F14").Value = 61.18887088966495
$ws.Range("J14").Value = 10.66421118627684

$ws.Range("B15").Value = 22.96064051363398
$ws.Range("C15").Value = 15.13014253788652
$ws.Range("D15").Value = 5.369215277711696
$ws.Range("E15").Value = 11.61027069457949
$ws.Range("F15").Value = 61.064642105046
$ws.Range("J15").Value = 10.66251332281504

$ws.Range("B16").Value = 22.68509423064429
$ws.Range("C16").Value = 14.89696124279176
$ws.Range("D16").Value = 5.376216780773007
$ws.Range("E16").Value = 11.62229442858698
$ws.Range("F16").Value = 60.35150187716742
$ws.Range("J16").Value = 10.65337657001777

$ws.Range("B17").Value = 22.51686088568152
$ws.Range("C17").Value = 14.75396572697551
$ws.Range("D17").Value = 5.380870652416911
$ws.Range("E17").Value = 11.63051484536439
$ws.Range("F17").Value = 59.91308748018208
$ws.Range("J17").Value = 10.6483006064824

$ws.Range("B18").Value = 22.42043396983891
$ws.Range("C18").Value = 14.67177601581209
$ws.Range("D18").Value = 5.383678328514785
$ws.Range("E18").Value = 11.63555282671348
$ws.Range("F18").Value = 59.66065362951078
$ws.Range("J18").Value = 10.64557546991097

$ws.Range("B19").Value = 22.38784712926908
$ws.Range("C19").Value = 14.64396152654531
$ws.Range("D19").Value = 5.384651353007973
$ws.Range("E19").Value = 11.63731176798901
$ws.Range("F19").Value = 59.57514437638874
$ws.Range("J19").Value = 10.64468617240001

$ws.Range("B20").Value = 22.53473570896235
$ws.Range("C20").Value = 14.76918275443278
$ws.Range("D20").Value = 5.380361672251246
$ws.Range("E20").Value = 11.62960769302747
$ws.Range("F20").Value = 59.95978694254449
$ws.Range("J20").Value = 10.64882082823021

$ws.Range("B21").Value = 23.03199744578615
$ws.Range("C21").Value = 15.19033177075189
$ws.Range("D21").Value = 5.367519266288555
$ws.Range("E21").Value = 11.60742775561847
$ws.Range("F21").Value = 61.24841766541023
$ws.Range("J21").Value = 10.66503584749331

$ws.Range("B22").Value = 23.35913058884635
$ws.Range("C22").Value = 15.46530141010081
$ws.Range("D22").Value = 5.360311146176254
$ws.Range("E22").Value = 11.59569979343234
$ws.Range("F22").Value = 62.08672698572722
$ws.Range("J22").Value = 10.67736483197305

$ws.Range("B23").Value = 23.18437788051479
$ws.Range("C23").Value = 15.31860700261044
$ws.Range("D23").Value = 5.364048868553092
$ws.Range("E23").Value = 11.60170555530154
$ws.Range("F23").Value = 61.63973136927676
$ws.Range("J23").Value = 10.67062621541353

$ws.Range("B24").Value = 22.52665358018873
$ws.Range("C24").Value = 14.76230306412296
$ws.Range("D24").Value = 5.380591370928778
$ws.Range("E24").Value = 11.63001684503508
$ws.Range("F24").Value = 59.93867529365558
$ws.Range("J24").Value = 10.64858503451429

$ws.Range("B25").Value = 21.83058227052966
$ws.Range("C25").Value = 14.16496131831078
$ws.Range("D25").Value = 5.403427808182553
$ws.Range("E25").Value = 11.67239310460168
$ws.Range("F25").Value = 58.09445005625859
$ws.Range("J25").Value = 10.6322408869079
